# "upgrade left table until javakheti"
#
# This workbook is one sheet of a batch of per-municipality "Number of
# stillbirths" tables. The edit:
#   1. Gives the sheet tab a real name ("Tetritskaro") instead of the
#      generic "1".
#   2. Masks out the Urban/Rural breakdown (rows "Urban" and "Rural") with
#      the confidentiality placeholder used elsewhere in the sheet, leaving
#      the "Total" row untouched.
#   3. Removes the blank spacer row that used to sit between the data table
#      and the footnote, so the footnote moves up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Real sheet name instead of the placeholder "1".
$ws.Name = "Tetritskaro"

# 2. Mask the "Urban" (row 6) and "Rural" (row 7) detail rows. The first
#    data column keeps the literal "..." marker already used for the
#    workbook's blank year column B, the remaining columns reuse the
#    existing "…" confidentiality marker used throughout the sheet.
$ws.Range("B6").Value = "..."
$ws.Range("C6:O6").Value = "…"

$ws.Range("B7").Value = "..."
$ws.Range("C7:O7").Value = "…"

# 3. Drop the now-empty separator row so the footnote ("Note: ...") shifts
#    up from row 9 to row 8.
$ws.Rows(8).Delete()
